$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.058.42'
$ws.Range("E2").Value = '  -0.80%  '
$ws.Range("D3").Value = '1.897.53'
$ws.Range("E3").Value = '  -1.81%  '
$ws.Range("E4").Value = '  +0.42%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7274'
$ws.Range("E5").Value = '  -7.30%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.32'
$ws.Range("E6").Value = '  -1.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.001'
$ws.Range("E7").Value = '  +0.49%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3093'
$ws.Range("E8").Value = '  -3.79%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '26.14'
$ws.Range("E9").Value = '  -6.41%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.06876'
$ws.Range("E10").Value = '  -2.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.7696'
$ws.Range("E11").Value = '  -1.53%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07939'
$ws.Range("E12").Value = '  -0.67%  '
$ws.Range("D13").Value = '1.894.21'
$ws.Range("E13").Value = '  -1.94%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.237'
$ws.Range("E14").Value = '  -2.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '90.79'
$ws.Range("E15").Value = '  -4.26%  '
$ws.Range("D16").Value = '30.057.38'
$ws.Range("E16").Value = '  -0.77%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.09'
$ws.Range("E17").Value = '  -3.04%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.770'
$ws.Range("E18").Value = '  -0.70%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007743'
$ws.Range("E19").Value = '  -3.19%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '237.02'
$ws.Range("E20").Value = '  -7.27%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.001'
$ws.Range("E21").Value = '  +0.41%  '
$ws.Range("D22").Value = '2.162.63'
$ws.Range("E22").Value = '  -1.20%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.47%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.924'
$ws.Range("E24").Value = '  +2.31%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.281'
$ws.Range("E25").Value = '  -3.14%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '165.49'
$ws.Range("E26").Value = '  +1.21%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.88'
$ws.Range("E27").Value = '  -1.08%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.1267'
$ws.Range("E28").Value = '  -7.43%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.016'
$ws.Range("E29").Value = '  -11.50%  '
$ws.Range("E30").Value = '  -1.21%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.534'
$ws.Range("E31").Value = '  +1.38%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.282'
$ws.Range("E32").Value = '  -3.16%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.062'
$ws.Range("E33").Value = '  -1.81%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.05081'
$ws.Range("E34").Value = '  -1.84%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.269'
$ws.Range("E35").Value = '  -1.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.7328'
$ws.Range("E36").Value = '  -2.49%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.738'
$ws.Range("E37").Value = '  -0.99%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01914'
$ws.Range("E38").Value = '  -2.50%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.771'
$ws.Range("E39").Value = '  -0.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.333'
$ws.Range("E40").Value = '  -1.23%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '74.48'
$ws.Range("E41").Value = '  -5.55%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.4413'
$ws.Range("E42").Value = '  -2.30%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.918'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9995'
$ws.Range("E44").Value = '  +0.23%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.8333'
$ws.Range("E45").Value = '  -0.01%  '
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '100.69'
$ws.Range("E46").Value = '  -0.13%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '7.572'
$ws.Range("E47").Value = '  +0.93%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.721'
$ws.Range("E48").Value = '  -0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '37.66'
$ws.Range("E49").Value = '  +1.21%  '
$ws.Range("B50").Value = 'RocketPoolETH'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range("D50").Value = '2.048.93'
$ws.Range("E50").Value = '  -1.82%  '
$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '934.45'
$ws.Range("E51").Value = '  -4.88%  '
